$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to be written/stored as literal text (matches the
    # original inlineStr cells), so numeric-looking strings such as
    # "306.66" are not silently reinterpreted as numbers by Excel.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "42.824.35"
Set-TextValue $ws.Range("E2") "  -0.49%  "
Set-TextValue $ws.Range("D3") "2.304.76"
Set-TextValue $ws.Range("E3") "  +0.22%  "
Set-TextValue $ws.Range("E4") "  +0.00%  "
Set-TextValue $ws.Range("D5") "306.66"
Set-TextValue $ws.Range("E5") "  +2.15%  "
Set-TextValue $ws.Range("D6") "96.26"
Set-TextValue $ws.Range("E6") "  -1.52%  "
Set-TextValue $ws.Range("D7") "0.508"
Set-TextValue $ws.Range("E7") "  -2.17%  "
Set-TextValue $ws.Range("D9") "0.502"
Set-TextValue $ws.Range("E9") "  -2.49%  "
Set-TextValue $ws.Range("D10") "35.30"
Set-TextValue $ws.Range("E10") "  -2.16%  "
Set-TextValue $ws.Range("D11") "0.0791"
Set-TextValue $ws.Range("E11") "  +0.16%  "
Set-TextValue $ws.Range("D12") "18.58"
Set-TextValue $ws.Range("E12") "  +5.05%  "
Set-TextValue $ws.Range("E13") "  +1.29%  "
Set-TextValue $ws.Range("D14") "6.77"
Set-TextValue $ws.Range("E14") "  -1.50%  "
Set-TextValue $ws.Range("D15") "2.666.47"
Set-TextValue $ws.Range("E15") "  +0.33%  "
Set-TextValue $ws.Range("D16") "2.301.07"
Set-TextValue $ws.Range("E16") "  -2.63%  "
Set-TextValue $ws.Range("D17") "0.782"
Set-TextValue $ws.Range("E17") "  -0.70%  "
Set-TextValue $ws.Range("D18") "42.760.91"
Set-TextValue $ws.Range("E18") "  -0.39%  "
Set-TextValue $ws.Range("D19") "13.09"
Set-TextValue $ws.Range("E19") "  +2.41%  "
Set-TextValue $ws.Range("D20") "0.0₃0897"
Set-TextValue $ws.Range("E20") "  -1.21%  "
Set-TextValue $ws.Range("E21") "  -1.16%  "
Set-TextValue $ws.Range("D22") "67.43"
Set-TextValue $ws.Range("E22") "  -2.14%  "
Set-TextValue $ws.Range("D23") "236.16"
Set-TextValue $ws.Range("E23") "  -0.67%  "
Set-TextValue $ws.Range("E24") "  -0.72%  "
Set-TextValue $ws.Range("E25") "  +0.88%  "
Set-TextValue $ws.Range("E26") "  -0.03%  "
Set-TextValue $ws.Range("E27") "  +0.25%  "
Set-TextValue $ws.Range("D28") "25.16"
Set-TextValue $ws.Range("E28") "  +0.84%  "
Set-TextValue $ws.Range("E29") "  +16.99%  "
Set-TextValue $ws.Range("D30") "166.51"
Set-TextValue $ws.Range("E30") "  +0.98%  "
Set-TextValue $ws.Range("D31") "9.06"
Set-TextValue $ws.Range("E31") "  -0.56%  "
Set-TextValue $ws.Range("D32") "33.08"
Set-TextValue $ws.Range("E32") "  +0.29%  "
Set-TextValue $ws.Range("E33") "  +0.07%  "
Set-TextValue $ws.Range("E34") "  +0.29%  "
Set-TextValue $ws.Range("B35") "Celestia"
Set-TextValue $ws.Range("C35") "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-TextValue $ws.Range("D35") "17.84"
Set-TextValue $ws.Range("E35") "  -0.04%  "
Set-TextValue $ws.Range("B36") "Filecoin"
Set-TextValue $ws.Range("C36") "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D36") "4.99"
Set-TextValue $ws.Range("E36") "  -2.08%  "
Set-TextValue $ws.Range("E37") "  -0.47%  "
Set-TextValue $ws.Range("D38") "0.0694"
Set-TextValue $ws.Range("E38") "  -0.53%  "
Set-TextValue $ws.Range("E39") "  -1.27%  "
Set-TextValue $ws.Range("E40") "  -1.03%  "
Set-TextValue $ws.Range("E41") "  -1.03%  "
Set-TextValue $ws.Range("E42") "  -2.53%  "
Set-TextValue $ws.Range("D43") "2.013.19"
Set-TextValue $ws.Range("E43") "  -0.12%  "
Set-TextValue $ws.Range("D44") "0.0279"
Set-TextValue $ws.Range("E44") "  -2.48%  "
Set-TextValue $ws.Range("D45") "18.37"
Set-TextValue $ws.Range("E45") "  +5.13%  "
Set-TextValue $ws.Range("E46") "  -3.02%  "
Set-TextValue $ws.Range("D47") "2.04"
Set-TextValue $ws.Range("E47") "  -7.11%  "
Set-TextValue $ws.Range("D48") "2.80"
Set-TextValue $ws.Range("E48") "  -0.77%  "
Set-TextValue $ws.Range("B49") "HuobiToken"
Set-TextValue $ws.Range("C49") "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue $ws.Range("D49") "2.87"
Set-TextValue $ws.Range("E49") "  +9.76%  "
Set-TextValue $ws.Range("B50") "MultiversX"
Set-TextValue $ws.Range("C50") "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
Set-TextValue $ws.Range("D50") "53.90"
Set-TextValue $ws.Range("E50") "  -0.33%  "
Set-TextValue $ws.Range("D51") "2.532.00"
Set-TextValue $ws.Range("E51") "  +0.20%  "
